# Creación de query para inserts en compromisos
#
# Adds a new data row (row 3) to the "Compromisos" sheet with sample
# values, formats the date cell, makes the header row (row 2) a bit
# taller, and simplifies the page setup (no forced fit-to-page scale).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 3) -------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Hola"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = "No se"
$ws.Range("F3").Value = "Tampoco"

# Date column: store as a real date (serial number 43341 = 2018-08-29)
# formatted with the built-in short-date number format (numFmtId 14).
$ws.Range("G3").Value = 43341
$ws.Range("G3").NumberFormat = "mm-dd-yy"

$ws.Range("H3").Value = "Si"
$ws.Range("I3").Value = "No"
$ws.Range("J3").Value = "Si"
$ws.Range("K3").Value = "Si"

# --- Header row (row 2) gets a bit taller ---------------------------------
$ws.Rows.Item(2).RowHeight = 30

# --- Page setup: rely on normal 100% zoom instead of fit-to-page ---------
$ws.PageSetup.Zoom = $true

Write-Host "Row 3 added and formatting updated."
